$d = $word.ActiveDocument

# The document contains exactly one table. One of its cells holds a
# two-token comparison rendered as "A3br1c1>A4br1c1" (three runs:
# "A3br1c1", ">", "A4br1c1"). The edit swaps the two tokens so the
# cell reads "A4br1c1>A3br1c1" instead (the ">" stays put).
$oldText = "A3br1c1>A4br1c1"
$newText = "A4br1c1>A3br1c1"

$targetCell = $null

foreach ($table in $d.Tables) {
    foreach ($row in $table.Rows) {
        foreach ($cell in $row.Cells) {
            $cellText = $cell.Range.Text.TrimEnd([char]0x07, [char]0x0D, [char]0x0A)
            if ($cellText -eq $oldText) {
                $targetCell = $cell
            }
        }
    }
}

if ($targetCell -eq $null) {
    throw "Could not locate a cell whose text is exactly '$oldText'"
}

$cellRange = $targetCell.Range
$textRange = $d.Range($cellRange.Start, $cellRange.Start + $oldText.Length)
$textRange.Text = $newText

Write-Output ("Updated cell text: " + $targetCell.Range.Text.TrimEnd([char]0x07, [char]0x0D, [char]0x0A))
